$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9 (year 2025) values per updated BIBI invoicing data
$ws.Range("B9").Value = 3973856.42
$ws.Range("C9").Value = 622826.24
$ws.Range("D9").Value = 4596682.66
$ws.Range("E9").Value = 13.54947221873263
$ws.Range("F9").Value = 86.45052778126737
$ws.Range("G9").Value = -39.80694945880388
$ws.Range("H9").Value = -28.23764587656619
$ws.Range("I9").Value = 39942
$ws.Range("J9").Value = 1713
$ws.Range("K9").Value = 41655
$ws.Range("L9").Value = 28842
$ws.Range("M9").Value = 159.3746154912974
$ws.Range("N9").Value = 8.808296094951862
